$wb = $excel.ActiveWorkbook

# Add the new "Tries" worksheet after the last existing sheet ("Trees")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Tries"

# Header row
$ws.Range("A1").Value = "Date Solved"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Algorithm"
$ws.Range("D1").Value = "Difficulty"
$ws.Range("E1").Value = "Solved First Time"
$ws.Range("F1").Value = "Video Help"
$ws.Range("G1").Value = "Revisit?"
$ws.Range("H1").Value = "Understand?"
$ws.Range("I1").Value = "Revisit Date #1"
$ws.Range("J1").Value = "Revisit Date #2"
$ws.Range("K1").Value = "Revisit Date #3"
$ws.Range("L1").Value = "Confidence Now"

# Row 2 - Implement Trie Prefix Tree
$ws.Range("A2").Value = "'09/30/25"
$ws.Range("B2").Value = "Implement Trie Prefix Tree"
$ws.Range("C2").Value = "Tries"
$ws.Range("D2").Value = "Medium"
$ws.Range("E2").Value = "Kinda"
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "Maybe"
$ws.Range("H2").Value = "Yes"

# Row 3 - Design Add and Search Words Data Structure
$ws.Range("A3").Value = "'"
$ws.Range("B3").Value = "Design Add and Search Words Data Structure"
$ws.Range("C3").Value = "Tries"
$ws.Range("D3").Value = "Medium"

# Row 4 - Word Search 2
$ws.Range("B4").Value = "Word Search 2"
$ws.Range("C4").Value = "Tries"
$ws.Range("D4").Value = "Hard"
